# Added Software Design Document
# Fills in the RACI chart rows for the Software Design Document sections
# (previously placeholder "…" rows) and removes the trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 keeps its "…" placeholder as-is; rows 10-19 get real content.
$ws.Range("A10").Value = "1.1 Problem Background"
$ws.Range("B10").Value = "R A"
$ws.Range("C10").Value = "C I"
$ws.Range("D10").Value = "C I"

$ws.Range("A11").Value = "1.2 System Overview"
$ws.Range("B11").Value = "R A"
$ws.Range("C11").Value = "C I"
$ws.Range("D11").Value = "C I"

$ws.Range("A12").Value = "1.3 Potential Benefits"
$ws.Range("B12").Value = "R A"
$ws.Range("C12").Value = "C I"
$ws.Range("D12").Value = "C I"

$ws.Range("A13").Value = "2.1 User Requirements"
$ws.Range("B13").Value = "C I"
$ws.Range("C13").Value = "R A"
$ws.Range("D13").Value = "C I"

$ws.Range("A14").Value = "2.2 Software Requirements"
$ws.Range("B14").Value = "C I"
$ws.Range("C14").Value = "R A"
$ws.Range("D14").Value = "C I"

$ws.Range("A15").Value = "2.3 Use Cases"
$ws.Range("B15").Value = "C I"
$ws.Range("C15").Value = "R A"
$ws.Range("D15").Value = "C I"

$ws.Range("A16").Value = "3.1 Software Design"
$ws.Range("B16").Value = "C I"
$ws.Range("C16").Value = "C I"
$ws.Range("D16").Value = "R A"

$ws.Range("A17").Value = "3.2 System Componenets"
$ws.Range("B17").Value = "C I"
$ws.Range("C17").Value = "C I"
$ws.Range("D17").Value = "R A"

$ws.Range("A18").Value = "4.1 Structural Design"
$ws.Range("B18").Value = "R A"
$ws.Range("C18").Value = "C I"
$ws.Range("D18").Value = "C I"

$ws.Range("A19").Value = "4.2 Visual Design"
$ws.Range("B19").Value = "R A"
$ws.Range("C19").Value = "C I"
$ws.Range("D19").Value = "C I"

# The sheet previously ended with an empty, bordered row 30; remove it so
# the used range shrinks back down to row 29.
$ws.Rows("30:30").Delete()

# Restore the active selection to match the new content region.
$ws.Range("E19").Select()
